# Commit: "New changes as of 04/15"
# Strip the trailing "LIMIT 100" clause from the three saved SQL queries
# (B2: participant/sample summary, B3: distinct sample/participant list,
#  B4: file/sample/participant listing), then move the active selection
# from C2 to C3 to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = @"
WITH Distinct_Samples AS (
    SELECT DISTINCT
        sp.participant_id,
        sp.study_participant_id,
        s.study_name,
        s.phs_accession,
        sp.gender,
        smp.sample_id
    FROM 
        df_participant sp
    JOIN 
        df_study s ON sp."study.phs_accession" = s.phs_accession
    JOIN 
        df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
    JOIN
        df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
    JOIN
        df_program p ON p.program_acronym = s."program.program_acronym"
    JOIN
        df_file f1 ON f1."sample.sample_id" = smp.sample_id
    JOIN
        df_genomic_info gi ON gi."file.file_id" = f1.file_id
    WHERE 
        s.phs_accession = 'phs001287'
),
Sample_Limit AS (
    SELECT
        participant_id,
        study_participant_id,
        study_name,
        phs_accession,
        gender,
        sample_id,
        ROW_NUMBER() OVER (PARTITION BY study_participant_id ORDER BY sample_id) as rn
    FROM 
        Distinct_Samples
)
SELECT
    participant_id AS "Participant ID", 
    study_name AS "Study Name",
    phs_accession AS Accession,
    gender AS Gender,
    -- Handle sample IDs with proper LIMIT and concatenation
    GROUP_CONCAT(
        CASE 
            WHEN rn <= 5 THEN sample_id 
            ELSE NULL  -- Ensures only 5 samples are considered
        END, ', '
    ) ||
    -- Add "..." if there are more than 5 samples for this participant
    CASE 
        WHEN MAX(rn) > 5 THEN ', ...' 
        ELSE '' 
    END AS Samples
FROM 
    Sample_Limit
GROUP BY
    participant_id, 
    study_name,
    phs_accession,
    gender
;
"@
# Writing the new (still very long, 62-line) text recalculates the
# wrapped-text autofit height; Excel caps row height at 409.5pt, so restore
# that cap explicitly (it was already at the cap before this edit too).
$ws.Rows.Item(2).RowHeight = 409.5

$ws.Range("B3").Value = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs001287'
ORDER BY 
    smp.sample_id ASC
;
"@

$ws.Range("B4").Value = @"
SELECT
    f1.file_name AS "File Name",
    s.study_name AS "Study Name",
    s.phs_accession AS "Accession",
    sp.participant_id AS "Participant Id",
    COALESCE((
        SELECT
            REPLACE(GROUP_CONCAT(CASE WHEN rn <= 5 THEN smp.sample_id ELSE NULL END, ', '), ', , ', ', ') ||
            CASE WHEN MAX(rn) > 5 THEN ', ...' ELSE '' END
        FROM (
            SELECT
                smp.sample_id,
                ROW_NUMBER() OVER (ORDER BY smp.sample_id) AS rn
            FROM df_sample smp
            WHERE smp."participant.study_participant_id" = sp.study_participant_id
        ) smp
    ), '') AS "Sample Id",
    f1.file_type AS "File Type",
    gi.library_strategy AS "Library Strategy"
FROM 
    df_study s
INNER JOIN 
    df_participant sp ON sp."study.phs_accession" = s.phs_accession
INNER JOIN  
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
INNER JOIN 
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
INNER JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
INNER JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
INNER JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
WHERE 
    s.phs_accession = 'phs001287' 
GROUP BY
    f1.file_name,
    s.study_name,
    s.phs_accession,
    sp.participant_id,
    f1.file_type,
    gi.library_strategy
ORDER BY 
    f1.file_name ASC;
"@

# Move the selected cell from C2 to C3, matching the saved cursor position.
$ws.Range("C3").Select()
